$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.574538
$ws.Range("H2").Value = 1.723614
$ws.Range("K2").Value = 2
$ws.Range("M2").Value = 0.6896435000000001
$ws.Range("N2").Value = 1.379287
$ws.Range("O2").Value = 0.03192461458046126
$ws.Range("P2").Value = 0.02377638465777991
$ws.Range("Q2").Value = 0.396226397203
$ws.Range("R2").Value = 2.377358383218
$ws.Range("S2").Value = 0.03192461458046126
$ws.Range("T2").Value = 0.02377638465777991

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.574538
$ws.Range("H3").Value = 1.723614
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.068283333333333
$ws.Range("N3").Value = 9.20485
$ws.Range("O3").Value = 0.1420353600669325
$ws.Range("P3").Value = 0.1586747749505109
$ws.Range("Q3").Value = 1.762845369766667
$ws.Range("R3").Value = 15.8656083279
$ws.Range("S3").Value = 0.1420353600669325
$ws.Range("T3").Value = 0.1586747749505109

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.574538
$ws.Range("H4").Value = 1.723614
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.304132666666667
$ws.Range("N4").Value = 15.912398
$ws.Range("O4").Value = 0.2455361227459804
$ws.Range("P4").Value = 0.2743006319030685
$ws.Range("Q4").Value = 3.047425774041333
$ws.Range("R4").Value = 27.426831966372
$ws.Range("S4").Value = 0.2455361227459804
$ws.Range("T4").Value = 0.2743006319030685

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.574538
$ws.Range("H5").Value = 1.723614
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 5.490996
$ws.Range("N5").Value = 16.472988
$ws.Range("O5").Value = 0.2541863019993003
$ws.Range("P5").Value = 0.2839641779781819
$ws.Range("Q5").Value = 3.154785859848
$ws.Range("R5").Value = 28.393072738632
$ws.Range("S5").Value = 0.2541863019993003
$ws.Range("T5").Value = 0.2839641779781819

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.574538
$ws.Range("H6").Value = 1.723614
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.9428843333333333
$ws.Range("N6").Value = 2.828653
$ws.Range("O6").Value = 0.04364750619069392
$ws.Range("P6").Value = 0.0487608030753448
$ws.Range("Q6").Value = 0.5417228791046667
$ws.Range("R6").Value = 4.875505911942
$ws.Range("S6").Value = 0.04364750619069392
$ws.Range("T6").Value = 0.0487608030753448

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.574538
$ws.Range("H7").Value = 1.723614
$ws.Range("K7").Value = 2
$ws.Range("M7").Value = 6.106310000000001
$ws.Range("N7").Value = 12.21262
$ws.Range("O7").Value = 0.2826700944166318
$ws.Range("P7").Value = 0.210523227435114
$ws.Range("Q7").Value = 3.50830713478
$ws.Range("R7").Value = 21.04984280868
$ws.Range("S7").Value = 0.2826700944166318
$ws.Range("T7").Value = 0.210523227435114

